# pTHg_SWeirR.xlsx update:
#  - Rename worksheet "6_pTHg" -> "7_pTHg" (folder moved from 6_pTHg to 7_pTHg / 7_rloadest)
#  - Correct cell C30 (row 30, "Q" value) from 11138 to 11134

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the new "7_pTHg" folder/stage naming
$ws.Name = "7_pTHg"

# Update the corrected data value in C30
$ws.Range("C30").Value = 11134
